$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get shuffled between rows (the rest of each row is
# identical across the whole table, so only these columns actually change).
$cols = @("D", "M", "N", "O", "P", "S")

# Mapping: row -> row whose OLD values become this row's NEW values.
$mapping = @{
    2  = 5
    3  = 6
    4  = 9
    5  = 2
    6  = 12
    7  = 11
    8  = 3
    9  = 10
    10 = 7
    11 = 13
    12 = 8
    13 = 4
}

# Snapshot the original values for the affected columns/rows before
# overwriting anything. Value2 is used (rather than Value) so that we get
# plain numeric values back instead of a wrapped Variant object.
$original = @{}
foreach ($row in 2..13) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowVals
}

# Apply the permutation.
foreach ($row in 2..13) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $original[$srcRow][$col]
    }
}
